$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Stash C5's original style (quote-prefixed variant) in a scratch cell so we
# can restore it after writing new text resets the cell's style.
$ws.Range("C5").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# Row 4 updates
$ws.Range("D4").Value = "item quantity added successfully"
$ws.Range("E4").Value = "Pass"

# Row 5 updates
$ws.Range("B5").Value = "Check cart updates (Following TC_03 Steps)"
$ws.Range("C5").Value = "Items added matches the cart content (Following TC_03 steps)"
$ws.Range("D5").Value = "Follow up cart successfully checked "
$ws.Range("E5").Value = "Pass"

# Restore C5's original style (lost when .Value was assigned above)
$ws.Range("H1").Copy()
$ws.Range("C5").PasteSpecial(-4122)

# New row 6 (TC_05) - start by copying row 5's row-level formatting (height,
# thickBot, per-cell styles) down one row, then overwrite with new values.
$ws.Range("A5:F5").Copy()
$ws.Range("A6:F6").PasteSpecial(-4122)

$ws.Range("A6").Value = "TC_05"
$ws.Range("B6").Value = "Remove item from cart (0 items in cart)"
$ws.Range("C6").Value = "Cart becomes empty after adding then removing one item"
$ws.Range("D6").Value = "item removal successful"
$ws.Range("E6").Value = "Pass"
$ws.Range("F6").Value = $null

# Fix C6's style too (writing the .Value above reset it to the non-prefixed
# variant; restore the style-4 variant from the stash)
$ws.Range("H1").Copy()
$ws.Range("C6").PasteSpecial(-4122)

# Clean up the scratch cell
$ws.Range("H1").Clear()

# Match the row height explicitly in case PasteSpecial didn't carry it
$ws.Rows.Item(6).RowHeight = 37.5

# Update selection to match target (single cell E5)
$ws.Range("E5").Select()
